$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix 1: "As at {{ doc_generated_date }}" used to be split across three runs
# ("As at {{ ", "doc_generated", "_date }}"). Collapse them back into a
# single run holding the whole phrase (all three already share identical
# run formatting, so a plain Find/Replace naturally merges them).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("As at {{ doc_generated_date }}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "As at {{ doc_generated_date }}", 2)

# ---------------------------------------------------------------------------
# Fix 2: the "Authorised Users of {{ mooring_name }} as at {{ issue_date }}"
# heading showed the wrong merge field (issue_date instead of
# doc_generated_date). Locate just the "issue_date" token inside that
# heading and swap it for "doc_generated_date".
# ---------------------------------------------------------------------------
$rngFull = $d.Content
$rngFull.Find.Execute("Authorised Users of {{ mooring_name }} as at {{ issue_date }}")
$headingStart = $rngFull.Start
$headingEnd = $rngFull.End

$rngToken = $d.Range($headingStart, $headingEnd)
$rngToken.Find.Execute("issue_date")
$tokenStart = $rngToken.Start
$tokenEnd = $rngToken.End

$rngToken = $d.Range($tokenStart, $tokenEnd)
$rngToken.Text = "doc_generated_date"
$rngToken.Font.Bold = $true
$rngToken.Font.Bold = $false
